$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.319712
$ws.Range("H2").Value = 3.959136
$ws.Range("Q2").Value = 0.188814275168
$ws.Range("R2").Value = 1.699328476512
